# Weekly fruit/vegetable price update: insert two new price rows
# (Especial / Primera quality, dated 2022-11-08) at the top of the
# data block, pushing the existing records down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 27 - all data that was
# in rows 27:41 shifts down to 29:43.
$ws.Rows("27:28").Insert()

# --- New row 27: Chirimoya "Especial" quality, week of 2022-11-08 ---
$ws.Range("A27").Value = 11
$ws.Range("B27").Value = "Vega Monumental Concepción"
$ws.Range("C27").Value = "Bíobío"
$ws.Range("D27").Value = 44873
$ws.Range("E27").Value = 8
$ws.Range("F27").Value = "Fruta"
$ws.Range("G27").Value = 100107
$ws.Range("H27").Value = "Otros"
$ws.Range("I27").Value = 100107002
$ws.Range("J27").Value = "Chirimoya"
$ws.Range("K27").Value = "Cultivar IV Región"
$ws.Range("L27").Value = "Especial"
$ws.Range("M27").Value = 50
$ws.Range("N27").Value = 26000
$ws.Range("O27").Value = 26000
$ws.Range("P27").Value = 26000
$ws.Range("Q27").Value = "$/bandeja 10 kilos"
$ws.Range("R27").Value = "Provincia de Limarí"
$ws.Range("S27").Value = 2600
$ws.Range("T27").Value = 10

# --- New row 28: Chirimoya "Primera" quality, week of 2022-11-08 ---
$ws.Range("A28").Value = 11
$ws.Range("B28").Value = "Vega Monumental Concepción"
$ws.Range("C28").Value = "Bíobío"
$ws.Range("D28").Value = 44873
$ws.Range("E28").Value = 8
$ws.Range("F28").Value = "Fruta"
$ws.Range("G28").Value = 100107
$ws.Range("H28").Value = "Otros"
$ws.Range("I28").Value = 100107002
$ws.Range("J28").Value = "Chirimoya"
$ws.Range("K28").Value = "Cultivar IV Región"
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 50
$ws.Range("N28").Value = 24000
$ws.Range("O28").Value = 24000
$ws.Range("P28").Value = 24000
$ws.Range("Q28").Value = "$/bandeja 10 kilos"
$ws.Range("R28").Value = "Provincia de Limarí"
$ws.Range("S28").Value = 2400
$ws.Range("T28").Value = 10
